$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.195.03"
$ws.Range("E2").Value = "  -0.68%  "
$ws.Range("D3").Value = "1.807.93"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "'312.65"
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("D7").Value = "'0.5128"
$ws.Range("E7").Value = "  -2.48%  "
$ws.Range("D8").Value = "'0.3950"
$ws.Range("E8").Value = "  +2.53%  "
$ws.Range("D9").Value = "'0.07791"
$ws.Range("E9").Value = "  -2.93%  "
$ws.Range("D10").Value = "'1.110"
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("D11").Value = "'41.05"
$ws.Range("E11").Value = "  -1.94%  "
$ws.Range("D12").Value = "'6.357"
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("D13").Value = "'0.9999"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("E14").Value = "  -2.15%  "
$ws.Range("D15").Value = "'7.338"
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("D16").Value = "1.803.71"
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("D17").Value = "'92.80"
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("D18").Value = "'0.00001077"
$ws.Range("E18").Value = "  -2.19%  "
$ws.Range("D19").Value = "'0.06564"
$ws.Range("E19").Value = "  -1.19%  "
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("D21").Value = "'17.29"
$ws.Range("E21").Value = "  -2.07%  "
$ws.Range("D22").Value = "'6.012"
$ws.Range("D23").Value = "28.252.83"
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("E24").Value = "  -2.12%  "
$ws.Range("D25").Value = "'2.216"
$ws.Range("E25").Value = "  -1.25%  "
$ws.Range("D26").Value = "'160.91"
$ws.Range("E26").Value = "  +1.00%  "
$ws.Range("D27").Value = "'2.459"
$ws.Range("E27").Value = "  +1.70%  "
$ws.Range("D28").Value = "'20.49"
$ws.Range("E28").Value = "  -1.72%  "
$ws.Range("D29").Value = "2.014.98"
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("D30").Value = "'127.88"
$ws.Range("E30").Value = "  +2.71%  "
$ws.Range("E31").Value = "  -1.11%  "
$ws.Range("E32").Value = "  -1.63%  "
$ws.Range("D33").Value = "'3.656"
$ws.Range("E33").Value = "  -0.59%  "
$ws.Range("D34").Value = "'5.569"
$ws.Range("E34").Value = "  -1.89%  "
$ws.Range("D35").Value = "'0.07137"
$ws.Range("E35").Value = "  -2.74%  "
$ws.Range("D36").Value = "'9.166"
$ws.Range("E36").Value = "  +4.66%  "
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").Value = "'0.2175"
$ws.Range("E38").Value = "  -1.13%  "
$ws.Range("D39").Value = "'5.041"
$ws.Range("E39").Value = "  -1.85%  "
$ws.Range("D40").Value = "'11.55"
$ws.Range("E40").Value = "  -5.57%  "
$ws.Range("E41").Value = "  -2.30%  "
$ws.Range("D42").Value = "'0.9999"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").Value = "'1.154"
$ws.Range("E43").Value = "  -2.47%  "
$ws.Range("D44").Value = "'13.13"
$ws.Range("E44").Value = "  -2.79%  "
$ws.Range("D45").Value = "'0.5959"
$ws.Range("E45").Value = "  -2.75%  "
$ws.Range("E46").Value = "  -5.54%  "
$ws.Range("D47").Value = "'3.737"
$ws.Range("E47").Value = "  -1.22%  "
$ws.Range("D48").Value = "'124.99"
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("D49").Value = "'1.208"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").Value = "'1.919"
$ws.Range("E50").Value = "  -3.18%  "
$ws.Range("D51").Value = "'0.06790"
$ws.Range("E51").Value = "  -1.57%  "
